# Update the RPAR_holdings workbook:
#  - Roll the "as of" date in the confidential disclaimer from 2021-04-23 to 2021-04-26
#  - Refresh the Weight (D) / Percent Change (E) figures for rows 2-15
#
# The sheet ships protected, so it must be unprotected before the cells can be
# written, then re-protected afterwards to restore the original protected state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# --- Disclaimer text (shared string) ---------------------------------------
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."

# --- Weight / Percent Change figures -----------------------------------------
$ws.Range("D2").Value = 0.05748042340189435
$ws.Range("E2").Value = 0.003637033285760261

$ws.Range("D3").Value = 0.02371090620509983
$ws.Range("E3").Value = 0.001562499999999911

$ws.Range("D4").Value = 0.03174178484415415
$ws.Range("E4").Value = 0.001313074470080666

$ws.Range("D5").Value = 0.03033130348173024
$ws.Range("E5").Value = 0.007459593866556258

$ws.Range("D6").Value = 0.0364158236928994
$ws.Range("E6").Value = 0.02507085240898199

$ws.Range("D7").Value = 0.01891792249512089
$ws.Range("E7").Value = 0.007603371333870346

$ws.Range("D8").Value = 0.004861595821656699
$ws.Range("E8").Value = 0.005103082261685987

$ws.Range("D9").Value = 0.006979633383476651
$ws.Range("E9").Value = -0.003981042654028433

$ws.Range("D10").Value = 0.070060964484377
$ws.Range("E10").Value = 0.002832861189801639

$ws.Range("D11").Value = 0.07018004827670173
$ws.Range("E11").Value = 0.002828054298642524

$ws.Range("D12").Value = 0.1482566751377804
$ws.Range("E12").Value = -0.001427959445951643

$ws.Range("D13").Value = 0.3863614100079382
$ws.Range("E13").Value = -0.001047943411055741

$ws.Range("D14").Value = 0.1147015087671704
$ws.Range("E14").Value = 0.002664728682170603

$ws.Range("E15").Value = 0.001653888633935674

$ws.Protect()
